# Update workbook "Ventas.xlsx": add new rows of data to the
# productos, clientes and ventas sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: productos
# ---------------------------------------------------------------
$productos = $wb.Worksheets.Item("productos")

# Row 2 (existing product, values updated)
$productos.Cells.Item(2,1).Value = "PR-5324400431X724"
$productos.Cells.Item(2,2).Value = "PRC-00001"
$productos.Cells.Item(2,3).Value = "Cereal"
$productos.Cells.Item(2,4).Value = 2
$productos.Cells.Item(2,5).Value = "Kellogs"
$productos.Cells.Item(2,6).Value = 26

# Row 3 (new product)
$productos.Cells.Item(3,1).Value = "PR-01062918706PY6"
$productos.Cells.Item(3,2).Value = "PRC-00002"
$productos.Cells.Item(3,3).Value = "Leche"
$productos.Cells.Item(3,4).Value = 19
$productos.Cells.Item(3,5).Value = "Lala"
$productos.Cells.Item(3,6).Value = 12

# Row 4 (new product)
$productos.Cells.Item(4,1).Value = "PR-2319481912D2GS"
$productos.Cells.Item(4,2).Value = "PRC-00003"
$productos.Cells.Item(4,3).Value = "Pan"
$productos.Cells.Item(4,4).Value = 22
$productos.Cells.Item(4,5).Value = "Bimbo"
$productos.Cells.Item(4,6).Value = 10

# ---------------------------------------------------------------
# Sheet 2: clientes
# ---------------------------------------------------------------
$clientes = $wb.Worksheets.Item("clientes")

# Row 2 (existing client, values updated)
$clientes.Cells.Item(2,1).Value = "CL-53599631145GRI"
$clientes.Cells.Item(2,2).Value = "CLC-00001"
$clientes.Cells.Item(2,3).Value = "Josefo"
$clientes.Cells.Item(2,4).Value = "Colonia Bonilla"
$clientes.Cells.Item(2,5).Value = "'33778899"
$clientes.Cells.Item(2,6).Value = "fer@gmail.com"
$clientes.Cells.Item(2,7).Value = "'2025-10-22"

# Row 3 (new client)
$clientes.Cells.Item(3,1).Value = "CL-0326828760RKBB"
$clientes.Cells.Item(3,2).Value = "CLC-00002"
$clientes.Cells.Item(3,3).Value = "Hermenegildo"
$clientes.Cells.Item(3,4).Value = "Puerto Barrios"
$clientes.Cells.Item(3,5).Value = "'21210909"
$clientes.Cells.Item(3,6).Value = "herqww@gmail.com"
$clientes.Cells.Item(3,7).Value = "'2025-10-22"

# Row 4 (new client)
$clientes.Cells.Item(4,1).Value = "CL-3437554150ZGOH"
$clientes.Cells.Item(4,2).Value = "CLC-00003"
$clientes.Cells.Item(4,3).Value = "Juancho"
$clientes.Cells.Item(4,4).Value = "Zimbabue"
$clientes.Cells.Item(4,5).Value = "'34778890"
$clientes.Cells.Item(4,6).Value = "juancho@gmail.com"
$clientes.Cells.Item(4,7).Value = "'2025-10-22"

# ---------------------------------------------------------------
# Sheet 3: ventas
# ---------------------------------------------------------------
$ventas = $wb.Worksheets.Item("ventas")

# Row 2
$ventas.Cells.Item(2,1).Value = 1
$ventas.Cells.Item(2,2).Value = "PRC-00001"
$ventas.Cells.Item(2,3).Value = "CLC-00001"
$ventas.Cells.Item(2,4).Value = 2
$ventas.Cells.Item(2,5).Value = 52
$ventas.Cells.Item(2,6).Value = "'2025-10-22 22:59:29"
$ventas.Cells.Item(2,7).Value = "anulada"

# Row 3
$ventas.Cells.Item(3,1).Value = 2
$ventas.Cells.Item(3,2).Value = "PRC-00002"
$ventas.Cells.Item(3,3).Value = "CLC-00002"
$ventas.Cells.Item(3,4).Value = 2
$ventas.Cells.Item(3,5).Value = 24
$ventas.Cells.Item(3,6).Value = "'2025-10-22 23:03:40"
$ventas.Cells.Item(3,7).Value = "anulada"

# Row 4
$ventas.Cells.Item(4,1).Value = 3
$ventas.Cells.Item(4,2).Value = "PRC-00003"
$ventas.Cells.Item(4,3).Value = "CLC-00001"
$ventas.Cells.Item(4,4).Value = 2
$ventas.Cells.Item(4,5).Value = 20
$ventas.Cells.Item(4,6).Value = "'2025-10-22 23:29:00"
$ventas.Cells.Item(4,7).Value = "anulada"

# Row 5
$ventas.Cells.Item(5,1).Value = 4
$ventas.Cells.Item(5,2).Value = "PRC-00002"
$ventas.Cells.Item(5,3).Value = "CLC-00001"
$ventas.Cells.Item(5,4).Value = 3
$ventas.Cells.Item(5,5).Value = 36
$ventas.Cells.Item(5,6).Value = "'2025-10-22 23:29:13"
$ventas.Cells.Item(5,7).Value = "activa"

# Row 6
$ventas.Cells.Item(6,1).Value = 5
$ventas.Cells.Item(6,2).Value = "PRC-00003"
$ventas.Cells.Item(6,3).Value = "CLC-00003"
$ventas.Cells.Item(6,4).Value = 1
$ventas.Cells.Item(6,5).Value = 10
$ventas.Cells.Item(6,6).Value = "'2025-10-22 23:35:20"
$ventas.Cells.Item(6,7).Value = "activa"
